# Restore theoretical-analysis source data (Posiciones / Velocidades / Aceleraciones)
# after fixing the video path (PV_5.mp4 -> PV_1.mp4) and related utils.py bugs.
# The underlying trajectory samples changed, so every Tiempo/X/Y (and their
# derived Vx/Vy, Ax/Ay) value is rewritten cell-by-cell, and the three extra
# trailing rows produced by the longer trajectory are appended.
$wb = $excel.ActiveWorkbook

# --- Sheet 1: Posiciones ---
$ws = $wb.Worksheets.Item("Posiciones")
$ws.Range("A2").Value = 1.483499999999998
$ws.Range("B2").Value = 4
$ws.Range("C2").Value = 130
$ws.Range("A3").Value = 1.515749999999998
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 130
$ws.Range("A4").Value = 1.547999999999998
$ws.Range("B4").Value = 11
$ws.Range("C4").Value = 116
$ws.Range("A5").Value = 1.580249999999998
$ws.Range("B5").Value = 11
$ws.Range("C5").Value = 116
$ws.Range("A6").Value = 1.612499999999998
$ws.Range("B6").Value = 19
$ws.Range("C6").Value = 107
$ws.Range("A7").Value = 1.644749999999998
$ws.Range("B7").Value = 19
$ws.Range("C7").Value = 107
$ws.Range("A8").Value = 1.676999999999998
$ws.Range("B8").Value = 32
$ws.Range("C8").Value = 104
$ws.Range("A9").Value = 1.709249999999998
$ws.Range("B9").Value = 32
$ws.Range("C9").Value = 104
$ws.Range("A10").Value = 1.741499999999998
$ws.Range("B10").Value = 46
$ws.Range("C10").Value = 107
$ws.Range("A11").Value = 1.773749999999997
$ws.Range("B11").Value = 46
$ws.Range("C11").Value = 107
$ws.Range("A12").Value = 1.805999999999997
$ws.Range("B12").Value = 61
$ws.Range("C12").Value = 116
$ws.Range("A13").Value = 1.838249999999997
$ws.Range("B13").Value = 61
$ws.Range("C13").Value = 116
$ws.Range("A14").Value = 1.870499999999997
$ws.Range("B14").Value = 75
$ws.Range("C14").Value = 132
$ws.Range("A15").Value = 1.902749999999997
$ws.Range("B15").Value = 75
$ws.Range("C15").Value = 132
$ws.Range("A16").Value = 1.934999999999997
$ws.Range("B16").Value = 89
$ws.Range("C16").Value = 153
$ws.Range("A17").Value = 1.967249999999997
$ws.Range("B17").Value = 89
$ws.Range("C17").Value = 153
$ws.Range("A18").Value = 1.999499999999997
$ws.Range("B18").Value = 104
$ws.Range("C18").Value = 179
$ws.Range("A19").Value = 2.031749999999997
$ws.Range("B19").Value = 104
$ws.Range("C19").Value = 179
$ws.Range("A20").Value = 2.063999999999997
$ws.Range("B20").Value = 117
$ws.Range("C20").Value = 211
$ws.Range("A21").Value = 2.096249999999996
$ws.Range("B21").Value = 117
$ws.Range("C21").Value = 211
$ws.Range("A22").Value = 2.128499999999996
$ws.Range("B22").Value = 117
$ws.Range("C22").Value = 211
$ws.Range("A23").Value = 2.160749999999996
$ws.Range("B23").Value = 130
$ws.Range("C23").Value = 248
$ws.Range("A24").Value = 2.192999999999996
$ws.Range("B24").Value = 130
$ws.Range("C24").Value = 248
$ws.Range("A25").Value = 2.225249999999996
$ws.Range("B25").Value = 142
$ws.Range("C25").Value = 290
$ws.Range("A26").Value = 2.257499999999996
$ws.Range("B26").Value = 142
$ws.Range("C26").Value = 290
$ws.Range("A27").Value = 2.289749999999996
$ws.Range("B27").Value = 154
$ws.Range("C27").Value = 337
$ws.Range("A28").Value = 2.321999999999996
$ws.Range("B28").Value = 154
$ws.Range("C28").Value = 337
$ws.Range("A29").Value = 2.354249999999996
$ws.Range("B29").Value = 165
$ws.Range("C29").Value = 389
$ws.Range("A30").Value = 2.386499999999995
$ws.Range("B30").Value = 165
$ws.Range("C30").Value = 389
$ws.Range("A31").Value = 2.418749999999995
$ws.Range("B31").Value = 176
$ws.Range("C31").Value = 437
$ws.Range("A32").Value = 2.450999999999995
$ws.Range("B32").Value = 176
$ws.Range("C32").Value = 437
$ws.Range("A33").Value = 2.483249999999995
$ws.Range("B33").Value = 183
$ws.Range("C33").Value = 402
$ws.Range("A34").Value = 2.515499999999995
$ws.Range("B34").Value = 183
$ws.Range("C34").Value = 402
$ws.Range("A35").Value = 2.547749999999995
$ws.Range("B35").Value = 183
$ws.Range("C35").Value = 402

# --- Sheet 2: Velocidades ---
$ws = $wb.Worksheets.Item("Velocidades")
$ws.Range("A2").Value = 1.515749999999998
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("A3").Value = 1.547999999999998
$ws.Range("B3").Value = 217.0542635658915
$ws.Range("C3").Value = -434.1085271317829
$ws.Range("A4").Value = 1.580249999999998
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("A5").Value = 1.612499999999998
$ws.Range("B5").Value = 248.062015503876
$ws.Range("C5").Value = -279.0697674418604
$ws.Range("A6").Value = 1.644749999999998
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("A7").Value = 1.676999999999998
$ws.Range("B7").Value = 403.1007751937984
$ws.Range("C7").Value = -93.02325581395348
$ws.Range("A8").Value = 1.709249999999998
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("A9").Value = 1.741499999999998
$ws.Range("B9").Value = 434.1085271317829
$ws.Range("C9").Value = 93.02325581395348
$ws.Range("A10").Value = 1.773749999999997
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 0
$ws.Range("A11").Value = 1.805999999999997
$ws.Range("B11").Value = 465.1162790697674
$ws.Range("C11").Value = 279.0697674418604
$ws.Range("A12").Value = 1.838249999999997
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 0
$ws.Range("A13").Value = 1.870499999999997
$ws.Range("B13").Value = 434.1085271317829
$ws.Range("C13").Value = 496.1240310077519
$ws.Range("A14").Value = 1.902749999999997
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 0
$ws.Range("A15").Value = 1.934999999999997
$ws.Range("B15").Value = 434.1085271317829
$ws.Range("C15").Value = 651.1627906976744
$ws.Range("A16").Value = 1.967249999999997
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 0
$ws.Range("A17").Value = 1.999499999999997
$ws.Range("B17").Value = 465.1162790697674
$ws.Range("C17").Value = 806.2015503875969
$ws.Range("A18").Value = 2.031749999999997
$ws.Range("B18").Value = 0
$ws.Range("C18").Value = 0
$ws.Range("A19").Value = 2.063999999999997
$ws.Range("B19").Value = 403.1007751937984
$ws.Range("C19").Value = 992.2480620155038
$ws.Range("A20").Value = 2.096249999999996
$ws.Range("B20").Value = 0
$ws.Range("C20").Value = 0
$ws.Range("A21").Value = 2.128499999999996
$ws.Range("B21").Value = 0
$ws.Range("C21").Value = 0
$ws.Range("A22").Value = 2.160749999999996
$ws.Range("B22").Value = 403.1007751937984
$ws.Range("C22").Value = 1147.286821705426
$ws.Range("A23").Value = 2.192999999999996
$ws.Range("B23").Value = 0
$ws.Range("C23").Value = 0
$ws.Range("A24").Value = 2.225249999999996
$ws.Range("B24").Value = 372.0930232558139
$ws.Range("C24").Value = 1302.325581395349
$ws.Range("A25").Value = 2.257499999999996
$ws.Range("B25").Value = 0
$ws.Range("C25").Value = 0
$ws.Range("A26").Value = 2.289749999999996
$ws.Range("B26").Value = 372.0930232558139
$ws.Range("C26").Value = 1457.364341085271
$ws.Range("A27").Value = 2.321999999999996
$ws.Range("B27").Value = 0
$ws.Range("C27").Value = 0
$ws.Range("A28").Value = 2.354249999999996
$ws.Range("B28").Value = 341.0852713178294
$ws.Range("C28").Value = 1612.403100775194
$ws.Range("A29").Value = 2.386499999999995
$ws.Range("B29").Value = 0
$ws.Range("C29").Value = 0
$ws.Range("A30").Value = 2.418749999999995
$ws.Range("B30").Value = 341.0852713178294
$ws.Range("C30").Value = 1488.372093023256
$ws.Range("A31").Value = 2.450999999999995
$ws.Range("B31").Value = 0
$ws.Range("C31").Value = 0
$ws.Range("A32").Value = 2.483249999999995
$ws.Range("B32").Value = 217.0542635658915
$ws.Range("C32").Value = -1085.271317829457
$ws.Range("A33").Value = 2.515499999999995
$ws.Range("B33").Value = 0
$ws.Range("C33").Value = 0
$ws.Range("A34").Value = 2.547749999999995
$ws.Range("B34").Value = 0
$ws.Range("C34").Value = 0

# --- Sheet 3: Aceleraciones ---
$ws = $wb.Worksheets.Item("Aceleraciones")
$ws.Range("A2").Value = 1.547999999999998
$ws.Range("B2").Value = 6730.364761733068
$ws.Range("C2").Value = -13460.72952346614
$ws.Range("A3").Value = 1.580249999999998
$ws.Range("B3").Value = -6730.364761733068
$ws.Range("C3").Value = 13460.72952346614
$ws.Range("A4").Value = 1.612499999999998
$ws.Range("B4").Value = 7691.845441980649
$ws.Range("C4").Value = -8653.326122228231
$ws.Range("A5").Value = 1.644749999999998
$ws.Range("B5").Value = -7691.845441980649
$ws.Range("C5").Value = 8653.326122228231
$ws.Range("A6").Value = 1.676999999999998
$ws.Range("B6").Value = 12499.24884321856
$ws.Range("C6").Value = -2884.442040742744
$ws.Range("A7").Value = 1.709249999999998
$ws.Range("B7").Value = -12499.24884321856
$ws.Range("C7").Value = 2884.442040742744
$ws.Range("A8").Value = 1.741499999999998
$ws.Range("B8").Value = 13460.72952346614
$ws.Range("C8").Value = 2884.442040742744
$ws.Range("A9").Value = 1.773749999999997
$ws.Range("B9").Value = -13460.72952346614
$ws.Range("C9").Value = -2884.442040742744
$ws.Range("A10").Value = 1.805999999999997
$ws.Range("B10").Value = 14422.21020371372
$ws.Range("C10").Value = 8653.326122228231
$ws.Range("A11").Value = 1.838249999999997
$ws.Range("B11").Value = -14422.21020371372
$ws.Range("C11").Value = -8653.326122228231
$ws.Range("A12").Value = 1.870499999999997
$ws.Range("B12").Value = 13460.72952346614
$ws.Range("C12").Value = 15383.6908839613
$ws.Range("A13").Value = 1.902749999999997
$ws.Range("B13").Value = -13460.72952346614
$ws.Range("C13").Value = -15383.6908839613
$ws.Range("A14").Value = 1.934999999999997
$ws.Range("B14").Value = 13460.72952346614
$ws.Range("C14").Value = 20191.0942851992
$ws.Range("A15").Value = 1.967249999999997
$ws.Range("B15").Value = -13460.72952346614
$ws.Range("C15").Value = -20191.0942851992
$ws.Range("A16").Value = 1.999499999999997
$ws.Range("B16").Value = 14422.21020371372
$ws.Range("C16").Value = 24998.49768643711
$ws.Range("A17").Value = 2.031749999999997
$ws.Range("B17").Value = -14422.21020371372
$ws.Range("C17").Value = -24998.49768643711
$ws.Range("A18").Value = 2.063999999999997
$ws.Range("B18").Value = 12499.24884321856
$ws.Range("C18").Value = 30767.3817679226
$ws.Range("A19").Value = 2.096249999999996
$ws.Range("B19").Value = -12499.24884321856
$ws.Range("C19").Value = -30767.3817679226
$ws.Range("A20").Value = 2.128499999999996
$ws.Range("B20").Value = 0
$ws.Range("C20").Value = 0
$ws.Range("A21").Value = 2.160749999999996
$ws.Range("B21").Value = 12499.24884321856
$ws.Range("C21").Value = 35574.7851691605
$ws.Range("A22").Value = 2.192999999999996
$ws.Range("B22").Value = -12499.24884321856
$ws.Range("C22").Value = -35574.7851691605
$ws.Range("A23").Value = 2.225249999999996
$ws.Range("B23").Value = 11537.76816297098
$ws.Range("C23").Value = 40382.18857039841
$ws.Range("A24").Value = 2.257499999999996
$ws.Range("B24").Value = -11537.76816297098
$ws.Range("C24").Value = -40382.18857039841
$ws.Range("A25").Value = 2.289749999999996
$ws.Range("B25").Value = 11537.76816297098
$ws.Range("C25").Value = 45189.59197163631
$ws.Range("A26").Value = 2.321999999999996
$ws.Range("B26").Value = -11537.76816297098
$ws.Range("C26").Value = -45189.59197163631
$ws.Range("A27").Value = 2.354249999999996
$ws.Range("B27").Value = 10576.28748272339
$ws.Range("C27").Value = 49996.99537287423
$ws.Range("A28").Value = 2.386499999999995
$ws.Range("B28").Value = -10576.28748272339
$ws.Range("C28").Value = -49996.99537287423
$ws.Range("A29").Value = 2.418749999999995
$ws.Range("B29").Value = 10576.28748272339
$ws.Range("C29").Value = 46151.0726518839
$ws.Range("A30").Value = 2.450999999999995
$ws.Range("B30").Value = -10576.28748272339
$ws.Range("C30").Value = -46151.0726518839
$ws.Range("A31").Value = 2.483249999999995
$ws.Range("B31").Value = 6730.364761733068
$ws.Range("C31").Value = -33651.82380866534
$ws.Range("A32").Value = 2.515499999999995
$ws.Range("B32").Value = -6730.364761733068
$ws.Range("C32").Value = 33651.82380866534
$ws.Range("A33").Value = 2.547749999999995
$ws.Range("B33").Value = 0
$ws.Range("C33").Value = 0

